$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "43.679.38"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.279.70"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  +0.31%  "

Set-TextValue $ws.Range("D5") "113.73"
$ws.Range("E5").Value = "  +9.96%  "

Set-TextValue $ws.Range("D6") "266.91"
$ws.Range("E6").Value = "  -1.56%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.11%  "

Set-TextValue $ws.Range("D9") "0.611"
$ws.Range("E9").Value = "  +0.71%  "

Set-TextValue $ws.Range("D10") "48.22"
$ws.Range("E10").Value = "  +5.19%  "

Set-TextValue $ws.Range("D11") "0.0934"
$ws.Range("E11").Value = "  -0.03%  "

Set-TextValue $ws.Range("D12") "8.82"
$ws.Range("E12").Value = "  +7.18%  "

Set-TextValue $ws.Range("D13") "0.107"
$ws.Range("E13").Value = "  -0.17%  "

$ws.Range("E14").Value = "  +0.76%  "

$ws.Range("D15").Value = "2.620.43"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").Value = "2.281.59"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").Value = "43.472.23"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("E20").Value = "  +11.55%  "

Set-TextValue $ws.Range("D21") "71.83"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("E22").Value = "  -4.26%  "

Set-TextValue $ws.Range("D23") "9.88"
$ws.Range("E23").Value = "  +6.33%  "

Set-TextValue $ws.Range("D24") "231.85"
$ws.Range("E24").Value = "  -0.73%  "

$ws.Range("E25").Value = "  -2.31%  "

$ws.Range("E26").Value = "  +0.00%  "

Set-TextValue $ws.Range("D27") "11.48"
$ws.Range("E27").Value = "  +1.58%  "

Set-TextValue $ws.Range("D28") "41.02"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("E30").Value = "  +1.44%  "

Set-TextValue $ws.Range("D31") "173.20"
$ws.Range("E31").Value = "  -2.53%  "

Set-TextValue $ws.Range("D32") "21.41"
$ws.Range("E32").Value = "  -1.92%  "

Set-TextValue $ws.Range("D33") "0.0912"
$ws.Range("E33").Value = "  +1.13%  "

Set-TextValue $ws.Range("D34") "5.62"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  -6.07%  "

Set-TextValue $ws.Range("D37") "0.0350"
$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  -5.88%  "

$ws.Range("E39").Value = "  +5.27%  "

Set-TextValue $ws.Range("D40") "14.30"
$ws.Range("E40").Value = "  +17.62%  "

Set-TextValue $ws.Range("D41") "74.36"
$ws.Range("E41").Value = "  +13.63%  "

Set-TextValue $ws.Range("D42") "2.42"
$ws.Range("E42").Value = "  +3.54%  "

$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("E44").Value = "  +15.99%  "

$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("E49").Value = "  +1.83%  "

Set-TextValue $ws.Range("D50") "101.38"
$ws.Range("E50").Value = "  +2.06%  "

$ws.Range("E51").Value = "  +3.09%  "
